# Ra_Stock_5 radium concentration correction.
# Any file with Ra_Stock_5 used had the wrong concentration for stock 5
# radium (Parameters!B6 "Stock Activity" + its error in C6) — fix it here.
# Every downstream formula (Bottle Results / Averaged Results) references
# Parameters!$B$6 / $C$6 directly, so they recalculate automatically.

$wb = $excel.ActiveWorkbook

$paramWs = $wb.Worksheets.Item("Parameters")
$paramWs.Activate()
$paramWs.Range("B6").Value = 1407
$paramWs.Range("C6").Value = 62
$paramWs.Range("C6").NumberFormat = "0.00E+00"
$paramWs.Range("B6:C6").Select()

$bottleWs = $wb.Worksheets.Item("Bottle Results")
$bottleWs.Activate()
$bottleWs.Range("Q23").Select()

$avgWs = $wb.Worksheets.Item("Averaged Results")
$avgWs.Activate()
$avgWs.Range("L17").Select()
